$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new row (row 18) with: eex_field_JSON = "url", machine_name = "field_link_api", is_dataset = FALSE
$ws.Range("A18").Value = "url"
$ws.Range("B18").Value = "field_link_api"
$ws.Range("C18").Value = $false

# Move the active selection to A19, matching the post-edit selection in the file
$ws.Range("A19").Select()
